$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Reorder rows: move the "Wave equation" row (currently row 8, no Nr.) to the
# bottom of the table (new row 13, also with no Nr.), shifting the rows below it
# (NSL, HUP, SLT, CL) up by one. ---
$ws.Rows.Item(8).Delete()

# Re-add the Wave equation entry at the end of the table (row 13), leaving row 12
# blank as a gap, and without a "Nr." number (consistent with its original entry).
$ws.Cells.Item(13, 2).Value = "Q193846"
$ws.Cells.Item(13, 3).Value = "WE"
$ws.Cells.Item(13, 4).Value = "Wave equation"

# --- Add the new "Ex. Equation" column (column E) with one example equation per
# row, taken from the first example row of each equation's own worksheet. ---
$ws.Range("E1").Value = "Ex. Equation"
$ws.Range("E1").Font.Bold = $true

$ws.Range("E2").Value = '\frac{1}{c^2} \frac{\partial^2 \psi}{\partial t^2} - \nabla^2 \psi + \left( \frac{m_0 c}{\hbar} \right)^2 \psi = 0'
$ws.Range("E3").Value = 'G_{\mu \nu} + \Lambda g_{\mu \nu} = \kappa T_{\mu \nu}'
$ws.Range("E4").Value = '\text{div} \vec{E} = 4 \pi \rho'
$ws.Range("E5").Value = 'i \hbar \frac{\partial}{\partial t} | \psi (t) \rangle = \hat{H} | \psi (t) \rangle'
$ws.Range("E6").Value = '(\nabla^2 − k^2) A = -f'
$ws.Range("E7").Value = '\nabla^4\varphi=0'
$ws.Range("E8").Value = '\vec{F} = \frac{d\vec{p}}{dt}'
$ws.Range("E9").Value = '\sigma_{x}\sigma_{p} \geq \frac{\hbar}{2}'
$ws.Range("E10").Value = '\oint \frac{\delta Q}{T}=0'
$ws.Range("E11").Value = '|F_1| = |F_2| = \\frac{|q_1 \\times q_2|}{r^2}'

$ws.Range("E2").Select()
